# Quarterly indexing esoteric bug-fix operation
# Updates the computed forecast-error metrics (ME, MAE, MSE, RMSE, SE, N)
# in rows 2-11 (columns B:G) of the active worksheet with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.5149680659622189, 0.5621490832682646, 0.398099301564984,  0.6309511086962158, 0.3773599919233752, 15),
    @(0.374392881790439,  0.3867115873339156, 0.2095230518030864, 0.4577368805362819, 0.2732908152567819, 14),
    @(0.3047392045302661, 0.3200881600767664, 0.1574586298430795, 0.3968105717380517, 0.2645285762523938, 13),
    @(0.4461589838474358, 0.4461589838474358, 0.2557785926581425, 0.5057455809576021, 0.248751253168957,  12),
    @(0.410700336148952,  0.4135187035930826, 0.2265423492358624, 0.4759646512461428, 0.2522981201580784, 11),
    @(0.3494740007469088, 0.3597227794264332, 0.1638948727701445, 0.4048393172236913, 0.2154138022300699, 10),
    @(0.3706545189354704, 0.3734090874777368, 0.1807682845776609, 0.4251685366741769, 0.2209218214474817, 9),
    @(0.3834281208488316, 0.385141425687041,  0.1960677877649145, 0.4427954242818171, 0.2367654991933658, 8),
    @(0.3533129028593847, 0.3533129028593847, 0.1644395682562605, 0.4055114896723155, 0.2149677985595936, 7),
    @(0.3494951052803065, 0.3494951052803065, 0.1636482358169425, 0.404534591619731,  0.2231629195060417, 6)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
